$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3714698851108551
$ws.Range("B1").Value = 0.2615610957145691
$ws.Range("C1").Value = 0.4057794809341431
$ws.Range("D1").Value = 4.472200393676758
$ws.Range("E1").Value = 2.317106008529663
